# New crime data collected — refresh weekly CompStat report numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) -----------------------------
# "Volume 30   Number  39" -> "...  40"
$ws.Range("A8").Characters(21, 2).Text = "40"

# "Report Covering the Week  9/25/2023  Through  10/1/2023"
#   -> "...10/2/2023  Through  10/8/2023"
$ws.Range("C9").Characters(27, 9).Text = "10/2/2023"
$ws.Range("C9").Characters(47, 9).Text = "10/8/2023"

# --- Row 14 (Murder) -------------------------------------------------------
$ws.Range("L14").Value = -80

# --- Row 15 (Rape): F/G/H switch from numbers to the "no data" text markers
$ws.Range("F15").Value = "'0"
$ws.Range("G15").Value = "'0"
$ws.Range("H15").Value = "'***.*"
$ws.Range("C15").Copy()
$ws.Range("F15:H15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 164
$ws.Range("J16").Value = 174
$ws.Range("K16").Value = -5.747126436781
$ws.Range("L16").Value = 10.810810810810
$ws.Range("M16").Value = -18
$ws.Range("N16").Value = -75.449101796407

# --- Row 17 (Fel. Assault) ---------------------------------------------------
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 278
$ws.Range("J17").Value = 278
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 7.751937984496
$ws.Range("M17").Value = 148.214285714286
$ws.Range("N17").Value = -13.931888544891

# --- Row 18 (Burglary) -------------------------------------------------------
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -55.555555555555
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 110
$ws.Range("J18").Value = 124
$ws.Range("K18").Value = -11.290322580645
$ws.Range("L18").Value = 50.684931506849
$ws.Range("M18").Value = -51.754385964912
$ws.Range("N18").Value = -87.654320987654

# --- Row 19 (Gr. Larceny) ----------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 27.272727272727
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 40
$ws.Range("I19").Value = 492
$ws.Range("J19").Value = 506
$ws.Range("K19").Value = -2.766798418972
$ws.Range("L19").Value = 30.851063829787
$ws.Range("M19").Value = 87.786259541984
$ws.Range("N19").Value = 10.313901345291

# --- Row 20 (G.L.A.) ----------------------------------------------------------
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 195
$ws.Range("J20").Value = 169
$ws.Range("K20").Value = 15.384615384615
$ws.Range("L20").Value = 39.285714285714
$ws.Range("M20").Value = -16.666666666666
$ws.Range("N20").Value = -92.212460063897

# --- Row 21 (TOTAL) ------------------------------------------------------------
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 7.692307692307
$ws.Range("F21").Value = 149
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = 15.503875968992
$ws.Range("I21").Value = 1255
$ws.Range("J21").Value = 1273
$ws.Range("K21").Value = -1.413982717989
$ws.Range("L21").Value = 24.134520276953
$ws.Range("M21").Value = 18.957345971564
$ws.Range("N21").Value = -74.256410256410

# --- Row 22 (Transit): D/E go from text markers to real numbers, F becomes the
#     new "no data" text marker -------------------------------------------------
$ws.Range("I22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D22").Value = 1

$ws.Range("L22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = -100

$ws.Range("F22").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -15

# --- Row 24 (Petit Larceny) -----------------------------------------------------
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 3.846153846153
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 994
$ws.Range("J24").Value = 1092
$ws.Range("K24").Value = -8.974358974358
$ws.Range("L24").Value = 43.434343434343
$ws.Range("M24").Value = 94.901960784313

# --- Row 25 (Misd. Assault) -----------------------------------------------------
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 62.5
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -10.256410256410
$ws.Range("I25").Value = 417
$ws.Range("J25").Value = 414
$ws.Range("K25").Value = 0.724637681159
$ws.Range("L25").Value = 14.246575342465
$ws.Range("M25").Value = 4.773869346733

# --- Row 26 (UCR Rape*): D/E go from text markers to real numbers ----------------
$ws.Range("I26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D26").Value = 2

$ws.Range("L26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E26").Value = -100

$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = -16.666666666666

# --- Row 27 (Other Sex Crimes): D/E go from text markers to real numbers --------
$ws.Range("I27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D27").Value = 4

$ws.Range("L27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = -100

$ws.Range("G27").Value = 7
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -19.565217391304
$ws.Range("L27").Value = 12.121212121212

# --- Row 28 (Shooting Vic.) ------------------------------------------------------
$ws.Range("L28").Value = -70.588235294117

# --- Row 29 (Shooting Inc.) ------------------------------------------------------
$ws.Range("L29").Value = -76.923076923076
